$d = $word.ActiveDocument

# --- Append four new paragraphs after the existing "No " paragraph --------
# Build them one at a time (InsertParagraphAfter + InsertAfter on a freshly
# collapsed end-of-document range) so each new paragraph gets its own <w:p>.

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("Contact developer:")

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("Name:")

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("Mob no:")

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

# Temporarily add trailing placeholder characters after "Email:" so the
# insertion point for the bookmark below is not the literal end of the
# document (adding a zero-length bookmark exactly at end-of-story lands it
# in the wrong place) -- the placeholder is stripped off again afterwards.
$r.InsertAfter("Email:XXX")

# --- Move the _GoBack bookmark from the old "No " paragraph to the end of
# the new "Email:" paragraph -------------------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$pc = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($pc)
$bmPos = $lastPara.Range.Start + 6   # right after "Email:", before "XXX"
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Remove the "XXX" placeholder now that the bookmark is anchored correctly.
$placeholder = $d.Range($bmPos, $lastPara.Range.End)
$placeholder.Delete()
